# DDAf_2023_Tableau_annexe_Tab06.xlsx - "Add files via upload" update
#
# 1) Country-name footnote-asterisk changes (shared strings), which in this
#    workbook mark "resource-rich" countries (see cell B2: "Pays (pays
#    riches en ressources ombrés)"). Rows whose country name carries a "*"
#    are shaded (light-blue fill); rows without are unshaded. The asterisk
#    toggles for three countries, so the corresponding row formatting
#    (column B..J fill/font) is toggled to match:
#      - "Soudan du Sud" (row 34, South Sudan) -> now marked resource-rich
#      - "Cabo Verde"     (row 48)             -> now marked resource-rich
#      - "Nigeria*"       (row 57)             -> no longer marked resource-rich
# 2) The regional/aggregate rows that roll up "resource-rich" vs.
#    "not resource-rich" countries are refreshed with their recomputed
#    figures (static values, there are no formulas in this sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1a. Soudan du Sud (row 34) becomes resource-rich: shade it like the
# other starred countries, e.g. Tchad* (row 17), then update its label.
$ws.Range("B17:J17").Copy()
$ws.Range("B34:J34").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B34").Value2 = "Soudan du Sud*"

# --- 1b. Cabo Verde (row 48) becomes resource-rich: it is already shaded
# like the other starred countries, so only the label needs the asterisk.
$ws.Range("B48").Value2 = "Cabo Verde*"

# --- 1c. Nigeria (row 57) is no longer resource-rich: remove the shading,
# matching an unshaded row such as République centrafricaine (row 16),
# then drop the asterisk from its label.
$ws.Range("B16:J16").Copy()
$ws.Range("B57:J57").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B57").Value2 = "Nigeria"

$excel.CutCopyMode = 0

# --- 2. Refresh the aggregate/region rows affected by the resource-rich
# reclassification above (values recomputed upstream; pasted verbatim).
$ws.Range("C69").Value2 = 81.5123514285714
$ws.Range("D69").Value2 = 80.5536928571429
$ws.Range("E69").Value2 = 82.5965985714286
$ws.Range("F69").Value2 = 0.97614714285714
$ws.Range("G69").Value2 = 69.8722942857143
$ws.Range("H69").Value2 = 64.0147557142857
$ws.Range("I69").Value2 = 76.2134342857143
$ws.Range("J69").Value2 = 0.83153285714286
$ws.Range("C77").Value2 = 98.8484845454546
$ws.Range("D77").Value2 = 98.9782054545455
$ws.Range("E77").Value2 = 98.7221036363637
$ws.Range("F77").Value2 = 1.00259090909091
$ws.Range("G77").Value2 = 94.3806354545455
$ws.Range("H77").Value2 = 93.5228781818182
$ws.Range("I77").Value2 = 95.2659927272727
$ws.Range("J77").Value2 = 0.98175454545455
$ws.Range("C80").Value2 = 73.43097375
$ws.Range("D80").Value2 = 70.40952875
$ws.Range("E80").Value2 = 76.6974275
$ws.Range("F80").Value2 = 0.8909
$ws.Range("G80").Value2 = 63.6907575
$ws.Range("H80").Value2 = 56.87362625
$ws.Range("I80").Value2 = 70.794945
$ws.Range("J80").Value2 = 0.7631425
$ws.Range("C82").Value2 = 79.576575952381
$ws.Range("D82").Value2 = 77.1294471428572
$ws.Range("E82").Value2 = 82.3013416666667
$ws.Range("F82").Value2 = 0.92154190476191
$ws.Range("G82").Value2 = 68.062125
$ws.Range("H82").Value2 = 61.997815952381
$ws.Range("I82").Value2 = 74.5787126190476
$ws.Range("J82").Value2 = 0.80435547619048
$ws.Range("E83").Value2 = 97.1821093939394
$ws.Range("C84").Value2 = 67.8614739130435
$ws.Range("D84").Value2 = 63.8516934782609
$ws.Range("E84").Value2 = 72.3573408695652
$ws.Range("F84").Value2 = 0.86088217391304
$ws.Range("G84").Value2 = 54.9493721739131
$ws.Range("H84").Value2 = 46.6706365217391
$ws.Range("I84").Value2 = 64.01821
$ws.Range("J84").Value2 = 0.70187304347826
$ws.Range("C86").Value2 = 85.2344604761905
$ws.Range("D86").Value2 = 83.3120814285715
$ws.Range("E86").Value2 = 87.2726966666667
$ws.Range("F86").Value2 = 0.94830857142857
$ws.Range("G86").Value2 = 74.2412414285714
$ws.Range("H86").Value2 = 68.710060952381
$ws.Range("I86").Value2 = 79.9204890476191
$ws.Range("J86").Value2 = 0.8477980952381
$ws.Range("C87").Value2 = 94.7722661538462
$ws.Range("D87").Value2 = 94.5784261538462
$ws.Range("E87").Value2 = 94.9509953846154
$ws.Range("F87").Value2 = 0.99491192307692
$ws.Range("G87").Value2 = 85.5596142307692
$ws.Range("H87").Value2 = 82.6056465384616
$ws.Range("I87").Value2 = 88.6114296153846
$ws.Range("J87").Value2 = 0.925055
$ws.Range("C89").Value2 = 98.8182516666667
$ws.Range("D89").Value2 = 98.9299233333334
$ws.Range("E89").Value2 = 98.7159305555556
$ws.Range("F89").Value2 = 1.00211583333333
$ws.Range("G89").Value2 = 95.9518619444444
$ws.Range("H89").Value2 = 95.3057347222222
$ws.Range("I89").Value2 = 96.6111422222222
$ws.Range("J89").Value2 = 0.98613166666667
$ws.Range("C90").Value2 = 99.1081014285714
$ws.Range("D90").Value2 = 99.3271028571429
$ws.Range("E90").Value2 = 98.9657100000001
$ws.Range("F90").Value2 = 1.00367
$ws.Range("G90").Value2 = 97.4816318181819
$ws.Range("H90").Value2 = 97.1373459090909
$ws.Range("I90").Value2 = 97.7446227272727
$ws.Range("J90").Value2 = 0.99384454545455
$ws.Range("H91").Value2 = 49.7276545161291
$ws.Range("C97").Value2 = 72.1731697058824
$ws.Range("D97").Value2 = 68.4322211764706
$ws.Range("E97").Value2 = 76.3165782352941
$ws.Range("F97").Value2 = 0.87812823529412
$ws.Range("G97").Value2 = 60.9303817647059
$ws.Range("H97").Value2 = 53.7724282352941
$ws.Range("I97").Value2 = 68.64605
$ws.Range("J97").Value2 = 0.75284735294118
$ws.Range("C98").Value2 = 90.12742625
$ws.Range("D98").Value2 = 88.889005625
$ws.Range("E98").Value2 = 91.412044375
$ws.Range("F98").Value2 = 0.965181875
$ws.Range("G98").Value2 = 79.61834125
$ws.Range("H98").Value2 = 75.60851
$ws.Range("I98").Value2 = 83.71196875
$ws.Range("J98").Value2 = 0.886136875
